$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 582
$ws.Cells.Item(582, 1).Value = 45189.94790576389
$ws.Cells.Item(582, 2).Value = "hshs0104746@naver.com"
$ws.Cells.Item(582, 3).Value = "광고홍보학과"
$ws.Cells.Item(582, 4).Value = 20232639
$ws.Cells.Item(582, 5).Value = "최희수"
$ws.Cells.Item(582, 6).Value = "74:26"
$ws.Cells.Item(582, 7).Value = 0.2
$ws.Cells.Item(582, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(582, 9).Value = "952만 명"
$ws.Cells.Item(582, 10).Value = 0.059
$ws.Cells.Item(582, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(582, 12).Value = "Red"
$ws.Cells.Item(582, 13).Value = "모름/무응답"

# Row 583
$ws.Cells.Item(583, 1).Value = 45189.97729207176
$ws.Cells.Item(583, 2).Value = "ehdus040127@naver.com"
$ws.Cells.Item(583, 3).Value = "사회복지학부"
$ws.Cells.Item(583, 4).Value = 20232307
$ws.Cells.Item(583, 5).Value = "김도연"
$ws.Cells.Item(583, 6).Value = "74:26"
$ws.Cells.Item(583, 7).Value = 0.2
$ws.Cells.Item(583, 8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(583, 9).Value = "166만 명"
$ws.Cells.Item(583, 10).Value = 0.002
$ws.Cells.Item(583, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(583, 12).Value = "Black"
$ws.Cells.Item(583, 14).Value = "모름/무응답"

# Row 584
$ws.Cells.Item(584, 1).Value = 45190.02308297454
$ws.Cells.Item(584, 2).Value = "hjkiubb@naver.com"
$ws.Cells.Item(584, 3).Value = "미디어스쿨"
$ws.Cells.Item(584, 4).Value = 20232538
$ws.Cells.Item(584, 5).Value = "박재은"
$ws.Cells.Item(584, 6).Value = "74:26"
$ws.Cells.Item(584, 7).Value = 0.2
$ws.Cells.Item(584, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(584, 9).Value = "952만 명"
$ws.Cells.Item(584, 10).Value = 0.059
$ws.Cells.Item(584, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(584, 12).Value = "Black"
$ws.Cells.Item(584, 14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"

# Row 585
$ws.Cells.Item(585, 1).Value = 45190.03219640046
$ws.Cells.Item(585, 2).Value = "smile001118@naver.com"
$ws.Cells.Item(585, 3).Value = "사회학과"
$ws.Cells.Item(585, 4).Value = 20222240
$ws.Cells.Item(585, 5).Value = "홍성준"
$ws.Cells.Item(585, 6).Value = "74:26"
$ws.Cells.Item(585, 7).Value = 0.2
$ws.Cells.Item(585, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(585, 9).Value = "952만 명"
$ws.Cells.Item(585, 10).Value = 0.059
$ws.Cells.Item(585, 11).Value = "상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다"
$ws.Cells.Item(585, 12).Value = "Red"
$ws.Cells.Item(585, 13).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"

# Row 586
$ws.Cells.Item(586, 1).Value = 45190.072204861106
$ws.Cells.Item(586, 2).Value = "poliku8630@naver.com"
$ws.Cells.Item(586, 3).Value = "컨텐츠 IT"
$ws.Cells.Item(586, 4).Value = 20205197
$ws.Cells.Item(586, 5).Value = "심지혁"
$ws.Cells.Item(586, 6).Value = "76:24"
$ws.Cells.Item(586, 7).Value = 0.2
$ws.Cells.Item(586, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(586, 9).Value = "952만 명"
$ws.Cells.Item(586, 10).Value = 0.059
$ws.Cells.Item(586, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(586, 12).Value = "Red"
$ws.Cells.Item(586, 13).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"

# Row 587
$ws.Cells.Item(587, 1).Value = 45190.07256223379
$ws.Cells.Item(587, 2).Value = "hyeeun7356@gmail.com"
$ws.Cells.Item(587, 3).Value = "식품영양학과"
$ws.Cells.Item(587, 4).Value = 20203824
$ws.Cells.Item(587, 5).Value = "유혜은"
$ws.Cells.Item(587, 6).Value = "74:26"
$ws.Cells.Item(587, 7).Value = 0.2
$ws.Cells.Item(587, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(587, 9).Value = "952만 명"
$ws.Cells.Item(587, 10).Value = 0.059
$ws.Cells.Item(587, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(587, 12).Value = "Black"
$ws.Cells.Item(587, 14).Value = "모름/무응답"

# Row 588
$ws.Cells.Item(588, 1).Value = 45190.14519309028
$ws.Cells.Item(588, 2).Value = "dnjsgmlwjd1020@naver.com"
$ws.Cells.Item(588, 3).Value = "인문학부"
$ws.Cells.Item(588, 4).Value = 20231057
$ws.Cells.Item(588, 5).Value = "원희정"
$ws.Cells.Item(588, 6).Value = "74:26"
$ws.Cells.Item(588, 7).Value = 0.2
$ws.Cells.Item(588, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(588, 9).Value = "952만 명"
$ws.Cells.Item(588, 10).Value = 0.059
$ws.Cells.Item(588, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(588, 12).Value = "Black"
$ws.Cells.Item(588, 14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"

# Row 589
$ws.Cells.Item(589, 1).Value = 45190.384077395836
$ws.Cells.Item(589, 2).Value = "minjeong7432@gmail.com"
$ws.Cells.Item(589, 3).Value = "간호학과"
$ws.Cells.Item(589, 4).Value = 20236217
$ws.Cells.Item(589, 5).Value = "김민정"
$ws.Cells.Item(589, 6).Value = "74:26"
$ws.Cells.Item(589, 7).Value = 0.2
$ws.Cells.Item(589, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(589, 9).Value = "952만 명"
$ws.Cells.Item(589, 10).Value = 0.059
$ws.Cells.Item(589, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(589, 12).Value = "Black"
$ws.Cells.Item(589, 14).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"

# Row 590
$ws.Cells.Item(590, 1).Value = 45190.41397086806
$ws.Cells.Item(590, 2).Value = "alsgk03@naver.com"
$ws.Cells.Item(590, 3).Value = "사회학과"
$ws.Cells.Item(590, 4).Value = 20222213
$ws.Cells.Item(590, 5).Value = "박민하"
$ws.Cells.Item(590, 6).Value = "74:26"
$ws.Cells.Item(590, 7).Value = 0.2
$ws.Cells.Item(590, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(590, 9).Value = "952만 명"
$ws.Cells.Item(590, 10).Value = 0.059
$ws.Cells.Item(590, 11).Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Cells.Item(590, 12).Value = "Black"
$ws.Cells.Item(590, 14).Value = "모름/무응답"

# Row 591
$ws.Cells.Item(591, 1).Value = 45190.52351451389
$ws.Cells.Item(591, 2).Value = "suani3176@gmail.com"
$ws.Cells.Item(591, 3).Value = "사회복지학부"
$ws.Cells.Item(591, 4).Value = 20232328
$ws.Cells.Item(591, 5).Value = "박수안"
$ws.Cells.Item(591, 6).Value = "75:25"
$ws.Cells.Item(591, 7).Value = 0.2
$ws.Cells.Item(591, 8).Value = "프랑스와 스웨덴의 국민부담률은 꾸준히 40%를 넘고 있다."
$ws.Cells.Item(591, 9).Value = "779만 명"
$ws.Cells.Item(591, 10).Value = 0.374
$ws.Cells.Item(591, 11).Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Cells.Item(591, 12).Value = "Black"
$ws.Cells.Item(591, 14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"

# Row 592
$ws.Cells.Item(592, 1).Value = 45190.588201342594
$ws.Cells.Item(592, 2).Value = "hkmcosmos1@gmail.com"
$ws.Cells.Item(592, 3).Value = "글로벌 비즈니스"
$ws.Cells.Item(592, 4).Value = 20226429
$ws.Cells.Item(592, 5).Value = "한기민"
$ws.Cells.Item(592, 6).Value = "77:23"
$ws.Cells.Item(592, 7).Value = 0.15
$ws.Cells.Item(592, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(592, 9).Value = "779만 명"
$ws.Cells.Item(592, 10).Value = 0.151
$ws.Cells.Item(592, 11).Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Cells.Item(592, 12).Value = "Black"
$ws.Cells.Item(592, 14).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"

# Row 593
$ws.Cells.Item(593, 1).Value = 45190.61006407408
$ws.Cells.Item(593, 2).Value = "yeon-jin22@naver.com"
$ws.Cells.Item(593, 3).Value = "데이터사이언스학부"
$ws.Cells.Item(593, 4).Value = 20233256
$ws.Cells.Item(593, 5).Value = "최연진"
$ws.Cells.Item(593, 6).Value = "74:26"
$ws.Cells.Item(593, 7).Value = 0.2
$ws.Cells.Item(593, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(593, 9).Value = "952만 명"
$ws.Cells.Item(593, 10).Value = 0.059
$ws.Cells.Item(593, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(593, 12).Value = "Black"
$ws.Cells.Item(593, 14).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"

# Row 594
$ws.Cells.Item(594, 1).Value = 45190.62784372685
$ws.Cells.Item(594, 2).Value = "bsw030409@naver.com"
$ws.Cells.Item(594, 3).Value = "철학과"
$ws.Cells.Item(594, 4).Value = 20221043
$ws.Cells.Item(594, 5).Value = "백승우"
$ws.Cells.Item(594, 6).Value = "74:26"
$ws.Cells.Item(594, 7).Value = 0.2
$ws.Cells.Item(594, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(594, 9).Value = "952만 명"
$ws.Cells.Item(594, 10).Value = 0.059
$ws.Cells.Item(594, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(594, 12).Value = "Black"
$ws.Cells.Item(594, 14).Value = "모름/무응답"

# Row 595
$ws.Cells.Item(595, 1).Value = 45190.63231763889
$ws.Cells.Item(595, 2).Value = "hyj4213@naver.com"
$ws.Cells.Item(595, 3).Value = "미디어스쿨"
$ws.Cells.Item(595, 4).Value = 20232590
$ws.Cells.Item(595, 5).Value = "함영준"
$ws.Cells.Item(595, 6).Value = "74:26"
$ws.Cells.Item(595, 7).Value = 0.2
$ws.Cells.Item(595, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(595, 9).Value = "952만 명"
$ws.Cells.Item(595, 10).Value = 0.059
$ws.Cells.Item(595, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(595, 12).Value = "Red"
$ws.Cells.Item(595, 13).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"

# Row 596
$ws.Cells.Item(596, 1).Value = 45190.69257645833
$ws.Cells.Item(596, 2).Value = "seo1020102p@naver.com"
$ws.Cells.Item(596, 3).Value = "미디어스쿨"
$ws.Cells.Item(596, 4).Value = 20232537
$ws.Cells.Item(596, 5).Value = "박재연"
$ws.Cells.Item(596, 6).Value = "74:26"
$ws.Cells.Item(596, 7).Value = 0.2
$ws.Cells.Item(596, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(596, 9).Value = "952만 명"
$ws.Cells.Item(596, 10).Value = 0.059
$ws.Cells.Item(596, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(596, 12).Value = "Red"
$ws.Cells.Item(596, 13).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"

# Row 597
$ws.Cells.Item(597, 1).Value = 45190.70098991898
$ws.Cells.Item(597, 2).Value = "kddong99@gmail.com"
$ws.Cells.Item(597, 3).Value = "빅데이터전공"
$ws.Cells.Item(597, 4).Value = 20181205
$ws.Cells.Item(597, 5).Value = "김상준"
$ws.Cells.Item(597, 6).Value = "76:24"
$ws.Cells.Item(597, 7).Value = 0.2
$ws.Cells.Item(597, 8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(597, 9).Value = "952만 명"
$ws.Cells.Item(597, 10).Value = 0.151
$ws.Cells.Item(597, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(597, 12).Value = "Red"
$ws.Cells.Item(597, 13).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"

# Row 598
$ws.Cells.Item(598, 1).Value = 45190.750975717594
$ws.Cells.Item(598, 2).Value = "psjj3840@gmail.com"
$ws.Cells.Item(598, 3).Value = "디지털미디어콘텐츠"
$ws.Cells.Item(598, 4).Value = 20215154
$ws.Cells.Item(598, 5).Value = "박서진"
$ws.Cells.Item(598, 6).Value = "74:26"
$ws.Cells.Item(598, 7).Value = 0.2
$ws.Cells.Item(598, 8).Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Cells.Item(598, 9).Value = "779만 명"
$ws.Cells.Item(598, 10).Value = 0.151
$ws.Cells.Item(598, 11).Value = "상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다"
$ws.Cells.Item(598, 12).Value = "Red"
$ws.Cells.Item(598, 13).Value = "모름/무응답"

# Row 599
$ws.Cells.Item(599, 1).Value = 45190.78272258102
$ws.Cells.Item(599, 2).Value = "ub030801@naver.com"
$ws.Cells.Item(599, 3).Value = "간호학과"
$ws.Cells.Item(599, 4).Value = 20226256
$ws.Cells.Item(599, 5).Value = "신유빈"
$ws.Cells.Item(599, 6).Value = "77:23"
$ws.Cells.Item(599, 7).Value = 0.15
$ws.Cells.Item(599, 8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(599, 9).Value = "779만 명"
$ws.Cells.Item(599, 10).Value = 0.374
$ws.Cells.Item(599, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(599, 12).Value = "Red"
$ws.Cells.Item(599, 13).Value = "모름/무응답"

# Row 600
$ws.Cells.Item(600, 1).Value = 45190.78992811343
$ws.Cells.Item(600, 2).Value = "yenaridia@naver.com"
$ws.Cells.Item(600, 3).Value = "경영학과"
$ws.Cells.Item(600, 4).Value = 20201634
$ws.Cells.Item(600, 5).Value = "최예나"
$ws.Cells.Item(600, 6).Value = "74:26"
$ws.Cells.Item(600, 7).Value = 0.2
$ws.Cells.Item(600, 8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(600, 9).Value = "952만 명"
$ws.Cells.Item(600, 10).Value = 0.059
$ws.Cells.Item(600, 11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(600, 12).Value = "Red"
$ws.Cells.Item(600, 13).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
# Format row 582 (template 578)
$ws.Range("A578:L578").Copy()
$ws.Range("A582:L582").PasteSpecial(-4122)
$ws.Range("M578").Copy()
$ws.Range("M582").PasteSpecial(-4122)
$ws.Rows.Item(582).RowHeight = 15.75

# Format row 583 (template 579)
$ws.Range("A579:L579").Copy()
$ws.Range("A583:L583").PasteSpecial(-4122)
$ws.Range("N579").Copy()
$ws.Range("N583").PasteSpecial(-4122)
$ws.Rows.Item(583).RowHeight = 15.75

# Format row 584 (template 579)
$ws.Range("A579:L579").Copy()
$ws.Range("A584:L584").PasteSpecial(-4122)
$ws.Range("N579").Copy()
$ws.Range("N584").PasteSpecial(-4122)
$ws.Rows.Item(584).RowHeight = 15.75

# Format row 585 (template 578)
$ws.Range("A578:L578").Copy()
$ws.Range("A585:L585").PasteSpecial(-4122)
$ws.Range("M578").Copy()
$ws.Range("M585").PasteSpecial(-4122)
$ws.Rows.Item(585).RowHeight = 15.75

# Format row 586 (template 578)
$ws.Range("A578:L578").Copy()
$ws.Range("A586:L586").PasteSpecial(-4122)
$ws.Range("M578").Copy()
$ws.Range("M586").PasteSpecial(-4122)
$ws.Rows.Item(586).RowHeight = 15.75

# Format row 587 (template 579)
$ws.Range("A579:L579").Copy()
$ws.Range("A587:L587").PasteSpecial(-4122)
$ws.Range("N579").Copy()
$ws.Range("N587").PasteSpecial(-4122)
$ws.Rows.Item(587).RowHeight = 15.75

# Format row 588 (template 579)
$ws.Range("A579:L579").Copy()
$ws.Range("A588:L588").PasteSpecial(-4122)
$ws.Range("N579").Copy()
$ws.Range("N588").PasteSpecial(-4122)
$ws.Rows.Item(588).RowHeight = 15.75

# Format row 589 (template 579)
$ws.Range("A579:L579").Copy()
$ws.Range("A589:L589").PasteSpecial(-4122)
$ws.Range("N579").Copy()
$ws.Range("N589").PasteSpecial(-4122)
$ws.Rows.Item(589).RowHeight = 15.75

# Format row 590 (template 579)
$ws.Range("A579:L579").Copy()
$ws.Range("A590:L590").PasteSpecial(-4122)
$ws.Range("N579").Copy()
$ws.Range("N590").PasteSpecial(-4122)
$ws.Rows.Item(590).RowHeight = 15.75

# Format row 591 (template 579)
$ws.Range("A579:L579").Copy()
$ws.Range("A591:L591").PasteSpecial(-4122)
$ws.Range("N579").Copy()
$ws.Range("N591").PasteSpecial(-4122)
$ws.Rows.Item(591).RowHeight = 15.75

# Format row 592 (template 579)
$ws.Range("A579:L579").Copy()
$ws.Range("A592:L592").PasteSpecial(-4122)
$ws.Range("N579").Copy()
$ws.Range("N592").PasteSpecial(-4122)
$ws.Rows.Item(592).RowHeight = 15.75

# Format row 593 (template 579)
$ws.Range("A579:L579").Copy()
$ws.Range("A593:L593").PasteSpecial(-4122)
$ws.Range("N579").Copy()
$ws.Range("N593").PasteSpecial(-4122)
$ws.Rows.Item(593).RowHeight = 15.75

# Format row 594 (template 579)
$ws.Range("A579:L579").Copy()
$ws.Range("A594:L594").PasteSpecial(-4122)
$ws.Range("N579").Copy()
$ws.Range("N594").PasteSpecial(-4122)
$ws.Rows.Item(594).RowHeight = 15.75

# Format row 595 (template 578)
$ws.Range("A578:L578").Copy()
$ws.Range("A595:L595").PasteSpecial(-4122)
$ws.Range("M578").Copy()
$ws.Range("M595").PasteSpecial(-4122)
$ws.Rows.Item(595).RowHeight = 15.75

# Format row 596 (template 578)
$ws.Range("A578:L578").Copy()
$ws.Range("A596:L596").PasteSpecial(-4122)
$ws.Range("M578").Copy()
$ws.Range("M596").PasteSpecial(-4122)
$ws.Rows.Item(596).RowHeight = 15.75

# Format row 597 (template 578)
$ws.Range("A578:L578").Copy()
$ws.Range("A597:L597").PasteSpecial(-4122)
$ws.Range("M578").Copy()
$ws.Range("M597").PasteSpecial(-4122)
$ws.Rows.Item(597).RowHeight = 15.75

# Format row 598 (template 578)
$ws.Range("A578:L578").Copy()
$ws.Range("A598:L598").PasteSpecial(-4122)
$ws.Range("M578").Copy()
$ws.Range("M598").PasteSpecial(-4122)
$ws.Rows.Item(598).RowHeight = 15.75

# Format row 599 (template 578)
$ws.Range("A578:L578").Copy()
$ws.Range("A599:L599").PasteSpecial(-4122)
$ws.Range("M578").Copy()
$ws.Range("M599").PasteSpecial(-4122)
$ws.Rows.Item(599).RowHeight = 15.75

# Format row 600 (template 578)
$ws.Range("A578:L578").Copy()
$ws.Range("A600:L600").PasteSpecial(-4122)
$ws.Range("M578").Copy()
$ws.Range("M600").PasteSpecial(-4122)
$ws.Rows.Item(600).RowHeight = 15.75
$excel.CutCopyMode = $false
$ws.Range("F578").Select()
